$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135, shifting rows 135:209 down to 136:210
$ws.Rows(135).Insert()

# Populate the new row 135 with the new record data
$ws.Range("A135").Value = 11
$ws.Range("B135").Value = "Vega Monumental Concepción"
$ws.Range("C135").Value = "Bíobío"
$ws.Range("D135").Value = 45097
$ws.Range("E135").Value = 8
$ws.Range("F135").Value = 100112021
$ws.Range("G135").Value = "Ají"
$ws.Range("H135").Value = "Inferno"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 25
$ws.Range("K135").Value = 12000
$ws.Range("L135").Value = 13000
$ws.Range("M135").Value = 12600
$ws.Range("N135").Value = "$/caja 15 kilos"
$ws.Range("O135").Value = "Región de Arica y Parinacota"
$ws.Range("P135").Value = 840
$ws.Range("Q135").Value = 15
$ws.Range("R135").Value = "Hortaliza"
